$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.082.62"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.133.71"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'589.00"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'147.76"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.126.93"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +12.67%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("D14").Value = "'37.57"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "3.652.52"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "63.929.16"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "'7.18"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "3.131.84"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'468.10"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").Value = "'0.733"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'13.29"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").Value = "'82.42"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'9.03"
$ws.Range("E27").Value = "  +9.48%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "'27.17"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "'0.109"
$ws.Range("D34").Value = "0.0₃0894"
$ws.Range("E34").Value = "  +11.39%  "
$ws.Range("E35").Value = "  +7.82%  "
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").Value = "'3.42"
$ws.Range("E37").Value = "  +12.26%  "
$ws.Range("D38").Value = "'6.11"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "'50.95"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "'455.68"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "2.897.00"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "'36.02"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").Value = "'125.42"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'24.82"
$ws.Range("E51").Value = "  +0.59%  "
